$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# TC38_Canine_Filter_Breed-Poodle.xlsx : "Fixed variables and query errors"
#
# The "CasesTab" Cypher query (row 2, column B) had two stray lines at its
# tail (a trailing comma on the `Response to Treatment` alias plus an extra
# `Cohort` column that the query never matched/needed). Those lines are
# removed here; the "SamplesTab" (row 3) and "FilesTab" (row 4) queries are
# left exactly as they were.
# ---------------------------------------------------------------------------

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Poodle']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Row heights settle a bit lower once the extra wrapped lines are gone.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 216
$ws.Rows.Item(4).RowHeight = 216

# The saved selection moves from B4 up to B2 (and the view no longer needs
# to keep row 4 pinned as the top-left cell).
$ws.Range("B2").Select() | Out-Null
